$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.19"
$ws.Range("E2").Value = "'2.15%"
$ws.Range("D3").Value = "'40.97"
$ws.Range("E3").Value = "'-0.09%"
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'-1.55%"
$ws.Range("D5").Value = "'0.07605"
$ws.Range("E5").Value = "'-0.73%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.701"
$ws.Range("E6").Value = "'4.00%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9310"
$ws.Range("E7").Value = "'1.73%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1200"
$ws.Range("E8").Value = "'-3.60%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1816"
$ws.Range("E9").Value = "'-0.39%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.08991"
$ws.Range("E10").Value = "'-1.06%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04139"
$ws.Range("E11").Value = "'-0.39%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.1052"
$ws.Range("E12").Value = "'0.29%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001287"
$ws.Range("E13").Value = "'2.36%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005842"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("B15").Value = "UpBots"
$ws.Range("C15").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D15").Value = "'0.007522"
$ws.Range("E15").Value = "'0.18%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.336"
$ws.Range("E16").Value = "'-0.27%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.327"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D19").Value = "'0.3347"
$ws.Range("E19").Value = "'0.36%"
$ws.Range("D20").Value = "'7.599"
$ws.Range("E20").Value = "'3.86%"
$ws.Range("D21").Value = "'0.1351"
$ws.Range("E21").Value = "'-2.37%"
$ws.Range("D22").Value = "'0.2837"
$ws.Range("E22").Value = "'-1.56%"
$ws.Range("D23").Value = "'0.03965"
$ws.Range("E23").Value = "'-2.76%"
$ws.Range("D24").Value = "'0.001281"
$ws.Range("E24").Value = "'1.34%"
$ws.Range("D25").Value = "'0.004081"
$ws.Range("E25").Value = "'-4.61%"
$ws.Range("D26").Value = "'0.0001350"
$ws.Range("E26").Value = "'6.16%"
$ws.Range("D38").Value = "'0.02421"
$ws.Range("E38").Value = "'-2.91%"
$ws.Range("D39").Value = "'0.05178"
$ws.Range("E39").Value = "'-2.84%"
$ws.Range("D40").Value = "'0.007705"
$ws.Range("E40").Value = "'-1.76%"
$ws.Range("D41").Value = "'0.1303"
$ws.Range("E41").Value = "'-0.59%"
$ws.Range("D42").Value = "'0.007588"
$ws.Range("E42").Value = "'10.10%"
$ws.Range("D43").Value = "'0.003300"
$ws.Range("E43").Value = "'72.52%"
$ws.Range("D44").Value = "'0.008475"
$ws.Range("E44").Value = "'10.80%"
$ws.Range("D45").Value = "'0.3392"
$ws.Range("E45").Value = "'10.91%"
$ws.Range("D46").Value = "'0.00006595"
$ws.Range("E46").Value = "'-1.80%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("E48").Value = "'58.17%"
$ws.Range("D49").Value = "'0.004200"
$ws.Range("E49").Value = "'35.27%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.15%"
